# Insert a new weekly record at the top of the date-ordered data block
# (row 21), pushing all existing records (rows 21-38) down by one row
# (to rows 22-39). The newly freed row 21 is then populated with the
# new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 21..38 down to 22..39
$ws.Rows.Item(21).Insert()

# Populate the new row 21 with the new weekly record
$ws.Cells.Item(21, 1).Value = 5
$ws.Cells.Item(21, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(21, 3).Value = "Maule"
$ws.Cells.Item(21, 4).Value = 44482
$ws.Cells.Item(21, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(21, 5).Value = 7
$ws.Cells.Item(21, 6).Value = 100112022
$ws.Cells.Item(21, 7).Value = "Arveja Verde"
$ws.Cells.Item(21, 8).Value = "Sin especificar"
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 200
$ws.Cells.Item(21, 11).Value = 22000
$ws.Cells.Item(21, 12).Value = 22000
$ws.Cells.Item(21, 13).Value = 22000
$ws.Cells.Item(21, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(21, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(21, 16).Value = 880
$ws.Cells.Item(21, 17).Value = 25
$ws.Cells.Item(21, 18).Value = "Hortaliza"
